# Update template's ontology terms
#
# 1. 4COM03_Metabolomics sheet: widen column B and mark the newly
#    introduced Term Source REF / Term Accession Number columns as
#    "user-specific" for rows 2-4.
# 2. SwateTemplateMetadata sheet: bump the recorded template Version from
#    1.1.6 to 1.1.7 (keeping the cell's existing number/quote-prefix
#    formatting intact).

$wb = $excel.ActiveWorkbook

$wsMain = $wb.Worksheets.Item("4COM03_Metabolomics")

# Widen column B (was bestFit at 31.42578125) to a fixed width of 37.
$wsMain.Columns.Item(2).ColumnWidth = 36.2

# Row 2: new "user-specific" ontology placeholder cells.
$wsMain.Range("C2").Value = "user-specific"
$wsMain.Range("D2").Value = "user-specific"
$wsMain.Range("F2").Value = "user-specific"
$wsMain.Range("G2").Value = "user-specific"
$wsMain.Range("I2").Value = "user-specific"
$wsMain.Range("J2").Value = "user-specific"
$wsMain.Range("L2").Value = "user-specific"
$wsMain.Range("M2").Value = "user-specific"

# Row 3.
$wsMain.Range("L3").Value = "user-specific"
$wsMain.Range("M3").Value = "user-specific"

# Row 4.
$wsMain.Range("L4").Value = "user-specific"
$wsMain.Range("M4").Value = "user-specific"

# Bump the template version number, preserving the cell's original style
# (which carries a quote-prefix format) by stashing/restoring it through a
# scratch cell around the value write.
$wsMeta = $wb.Worksheets.Item("SwateTemplateMetadata")
$scratch = $wsMeta.Range("Z1")
$versionCell = $wsMeta.Range("B3")

$versionCell.Copy($scratch)
$versionCell.Value = "1.1.7"
$scratch.Copy()
$versionCell.PasteSpecial(-4122)
$scratch.Clear()
